$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D (Categoria) from 11 to 14.
# Excel's ColumnWidth setter snaps to the pixel grid of the Normal-style
# font (adds ~0.8333 chars for integer inputs), so request 13.15 to land
# exactly on a stored width of 14 after that internal rounding.
$ws.Columns.Item(4).ColumnWidth = 13.15

# Update category label for "Pão francês"
$ws.Range("D2").Value = "Pães e doces"

# Update quantities (Quantidade column)
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 8
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 5
